$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1
$ws.Range("A1").Value = 45436

# Re-touch the A30:D30 merge so it is re-registered after A1:E1
$ws.Range("A30:D30").UnMerge()
$ws.Range("A30:D30").Merge()

# Update the "altas" price list (D23:D28)
$ws.Range("D23").Value = 11050
$ws.Range("D24").Value = 13650
$ws.Range("D25").Value = 19500
$ws.Range("D26").Value = 28600
$ws.Range("D27").Value = 52000
$ws.Range("D28").Value = 71500

# Update the "bajas" price list (D36:D37)
$ws.Range("D36").Value = 5920
$ws.Range("D37").Value = 8730
